$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing sheet to "Added Books" and clear its old content ---
$addedWs = $wb.Worksheets.Item(1)
$addedWs.Name = "Added Books"
$addedWs.Cells.ClearContents()

# --- Step 2: add a new sheet right after it, named "Changed Books" ---
$changedWs = $wb.Worksheets.Add($null, $addedWs, 1, 1)
$changedWs.Name = "Changed Books"

# --- Populate "Added Books" sheet ---
$addedWs.Range("A1").Value = "Name"
$addedWs.Range("B1").Value = "Link"
$addedWs.Range("C1").Value = "Price (vnd)"
$addedWs.Range("D1").Value = "Discount (%)"
$addedWs.Range("E1").Value = "Sold"
$addedWs.Range("F1").Value = "Rating"
$addedWs.Range("G1").Value = "Publisher"
$addedWs.Range("H1").Value = "Manufacturer"
$addedWs.Range("I1").Value = "Authors"
$addedWs.Range("J1").Value = "Other_sellers"
$addedWs.Range("A1:J1").Font.Bold = $true
$addedWs.Range("A1:J1").Borders.LineStyle = 1
$addedWs.Range("A1:J1").HorizontalAlignment = -4108
$addedWs.Range("A1:J1").VerticalAlignment = -4160

$addedWs.Range("A2").Value = "Bản Đồ (New Edition)"
$addedWs.Range("B2").Value = "https://tiki.vn/product-p50685547.html?spid=50685549"
$addedWs.Range("C2").Value = 224250
$addedWs.Range("D2").Value = 35
$addedWs.Range("E2").Value = 5935
$addedWs.Range("F2").Value = 4.8
$addedWs.Range("G2").Value = "['Nhã Nam']"
$addedWs.Range("H2").Value = "['Nhà Xuất Bản Lao Động']"
$addedWs.Range("I2").Value = "['Aleksandra Mizielińska', 'Daniel Mizieliński']"
$addedWs.Range("J2").Value = "[{'name': 'nha sach nguyet linh', 'price': 325000, 'link': 'https://tiki.vn/cua-hang/nha-sach-nguyet-linh'}, {'name': 'AHABOOKS', 'price': 310500, 'link': 'https://tiki.vn/cua-hang/ahabooks'}, {'name': 'Nhà Sách Trẻ Online', 'price': 344310, 'link': 'https://tiki.vn/cua-hang/nha-sach-tre-online'}, {'name': 'Phương Đông Books', 'price': 299500, 'link': 'https://tiki.vn/cua-hang/phuong-dong-books'}, {'name': 'HaAnBooks', 'price': 345000, 'link': 'https://tiki.vn/cua-hang/haanbooks'}, {'name': 'Nhà Sách Vĩnh Thụy', 'price': 311000, 'link': 'https://tiki.vn/cua-hang/nha-sach-vinh-thuy'}, {'name': 'Alpha Books Official', 'price': 276000, 'link': 'https://tiki.vn/cua-hang/alphabooks-official'}, {'name': 'SÁCH ĐẠI NAM', 'price': 340000, 'link': 'https://tiki.vn/cua-hang/sach-dai-nam'}, {'name': 'Sống Official', 'price': 276000, 'link': 'https://tiki.vn/cua-hang/abcbooks'}, {'name': 'Việt Thư Books', 'price': 327750, 'link': 'https://tiki.vn/cua-hang/viet-thu-books'}]"
$addedWs.Range("A2").Font.Bold = $true
$addedWs.Range("A2").Borders.LineStyle = 1
$addedWs.Range("A2").HorizontalAlignment = -4108
$addedWs.Range("A2").VerticalAlignment = -4160

$addedWs.Range("A3").Value = "Thiên Tài Bên Trái, Kẻ Điên Bên Phải (Tái Bản) (New Edition)"
$addedWs.Range("B3").Value = "https://tiki.vn/product-p109017985.html?spid=109017987"
$addedWs.Range("C3").Value = 110980
$addedWs.Range("D3").Value = 38
$addedWs.Range("E3").Value = 25756
$addedWs.Range("F3").Value = 4.8
$addedWs.Range("G3").Value = "['Vibooks']"
$addedWs.Range("H3").Value = "['Nhà Xuất Bản Thế Giới']"
$addedWs.Range("I3").Value = "['Cao Minh']"
$addedWs.Range("J3").Value = "[{'name': 'Skybooks Official Store', 'price': 161100, 'link': 'https://tiki.vn/cua-hang/skybooks-official-store'}, {'name': 'Tazano Official Store', 'price': 143000, 'link': 'https://tiki.vn/cua-hang/sachbanchay'}, {'name': 'AHABOOKS', 'price': 166351, 'link': 'https://tiki.vn/cua-hang/ahabooks'}, {'name': 'VBooks', 'price': 179000, 'link': 'https://tiki.vn/cua-hang/vbooks'}, {'name': 'Minhhabooks', 'price': 126000, 'link': 'https://tiki.vn/cua-hang/minhhabooks'}, {'name': 'Bamboo Books', 'price': 161100, 'link': 'https://tiki.vn/cua-hang/bamboo-books'}, {'name': 'Omega Plus Books', 'price': 143000, 'link': 'https://tiki.vn/cua-hang/omega-plus-books'}, {'name': 'Times Books', 'price': 151850, 'link': 'https://tiki.vn/cua-hang/times-books'}]"
$addedWs.Range("A3").Font.Bold = $true
$addedWs.Range("A3").Borders.LineStyle = 1
$addedWs.Range("A3").HorizontalAlignment = -4108
$addedWs.Range("A3").VerticalAlignment = -4160

$addedWs.Range("A4").Value = "Yêu Những Ngày Nắng Chẳng Ghét Những Ngày Mưa (New Edition)"
$addedWs.Range("B4").Value = "https://tiki.vn/product-p190861557.html?spid=190861559"
$addedWs.Range("C4").Value = 80520
$addedWs.Range("D4").Value = 39
$addedWs.Range("E4").Value = 5318
$addedWs.Range("F4").Value = 5
$addedWs.Range("G4").Value = "['Skybooks']"
$addedWs.Range("H4").Value = "['Nhà Xuất Bản Phụ Nữ']"
$addedWs.Range("I4").Value = "['Kulzsc']"
$addedWs.Range("J4").Value = "[{'name': 'AHABOOKS', 'price': 117500, 'link': 'https://tiki.vn/cua-hang/ahabooks'}, {'name': 'Nhà Sách Vĩnh Thụy', 'price': 132000, 'link': 'https://tiki.vn/cua-hang/nha-sach-vinh-thuy'}, {'name': 'Skybooks Official Store', 'price': 118800, 'link': 'https://tiki.vn/cua-hang/skybooks-official-store'}, {'name': 'Times Books', 'price': 111900, 'link': 'https://tiki.vn/cua-hang/times-books'}, {'name': 'SÁCH ĐẠI NAM', 'price': 132000, 'link': 'https://tiki.vn/cua-hang/sach-dai-nam'}]"
$addedWs.Range("A4").Font.Bold = $true
$addedWs.Range("A4").Borders.LineStyle = 1
$addedWs.Range("A4").HorizontalAlignment = -4108
$addedWs.Range("A4").VerticalAlignment = -4160

$addedWs.Range("A5").Value = "Càng Kỷ Luật, Càng Tự Do (New Edition)"
$addedWs.Range("B5").Value = "https://tiki.vn/product-p68585576.html?spid=68585577"
$addedWs.Range("C5").Value = 65400
$addedWs.Range("D5").Value = 40
$addedWs.Range("E5").Value = 21225
$addedWs.Range("F5").Value = 4.8
$addedWs.Range("G5").Value = "['Bloom Books']"
$addedWs.Range("H5").Value = "['Nhà Xuất Bản Thế Giới']"
$addedWs.Range("I5").Value = "['Ca Tây']"
$addedWs.Range("J5").Value = "[{'name': 'Nhà sách Fahasa', 'price': 96000, 'link': 'https://tiki.vn/cua-hang/nha-sach-fahasa'}, {'name': 'Skybooks Official Store', 'price': 98100, 'link': 'https://tiki.vn/cua-hang/skybooks-official-store'}, {'name': 'Nhà Sách Trẻ Online', 'price': 92050, 'link': 'https://tiki.vn/cua-hang/nha-sach-tre-online'}, {'name': 'Phương Đông Books', 'price': 99500, 'link': 'https://tiki.vn/cua-hang/phuong-dong-books'}, {'name': 'Times Books', 'price': 92350, 'link': 'https://tiki.vn/cua-hang/times-books'}, {'name': 'Omega Plus Books', 'price': 87000, 'link': 'https://tiki.vn/cua-hang/omega-plus-books'}, {'name': 'Nhà Sách Vĩnh Thụy', 'price': 98000, 'link': 'https://tiki.vn/cua-hang/nha-sach-vinh-thuy'}, {'name': 'Alpha Books Official', 'price': 87000, 'link': 'https://tiki.vn/cua-hang/alphabooks-official'}, {'name': 'VBooks', 'price': 109000, 'link': 'https://tiki.vn/cua-hang/vbooks'}, {'name': 'NHBook', 'price': 93000, 'link': 'https://tiki.vn/cua-hang/nhbook'}, {'name': 'Việt Thư Books', 'price': 103550, 'link': 'https://tiki.vn/cua-hang/viet-thu-books'}, {'name': 'info book', 'price': 98000, 'link': 'https://tiki.vn/cua-hang/info-book'}]"
$addedWs.Range("A5").Font.Bold = $true
$addedWs.Range("A5").Borders.LineStyle = 1
$addedWs.Range("A5").HorizontalAlignment = -4108
$addedWs.Range("A5").VerticalAlignment = -4160

$addedWs.Range("A6").Value = "Dear, Darling (New Edition)"
$addedWs.Range("B6").Value = "https://tiki.vn/product-p174444163.html?spid=174444165"
$addedWs.Range("C6").Value = 58960
$addedWs.Range("D6").Value = 33
$addedWs.Range("E6").Value = 6084
$addedWs.Range("F6").Value = 4.8
$addedWs.Range("G6").Value = "['Skybooks']"
$addedWs.Range("H6").Value = "['Nhà Xuất Bản Phụ Nữ Việt Nam']"
$addedWs.Range("I6").Value = "['Hiên']"
$addedWs.Range("J6").Value = "[{'name': 'Skybooks Official Store', 'price': 79200, 'link': 'https://tiki.vn/cua-hang/skybooks-official-store'}, {'name': 'NewShop Official', 'price': 70400, 'link': 'https://tiki.vn/cua-hang/newshopvn'}, {'name': 'Tazano Official Store', 'price': 68000, 'link': 'https://tiki.vn/cua-hang/sachbanchay'}, {'name': 'Phương Đông Books', 'price': 88000, 'link': 'https://tiki.vn/cua-hang/phuong-dong-books'}, {'name': 'AHABOOKS', 'price': 88000, 'link': 'https://tiki.vn/cua-hang/ahabooks'}, {'name': 'Nhà Sách Vĩnh Thụy', 'price': 88000, 'link': 'https://tiki.vn/cua-hang/nha-sach-vinh-thuy'}, {'name': 'info book', 'price': 88000, 'link': 'https://tiki.vn/cua-hang/info-book'}, {'name': 'ETS Books', 'price': 66000, 'link': 'https://tiki.vn/cua-hang/ets-books'}, {'name': 'Nhà Sách Trẻ Online', 'price': 74200, 'link': 'https://tiki.vn/cua-hang/nha-sach-tre-online'}, {'name': 'Times Books', 'price': 74500, 'link': 'https://tiki.vn/cua-hang/times-books'}]"
$addedWs.Range("A6").Font.Bold = $true
$addedWs.Range("A6").Borders.LineStyle = 1
$addedWs.Range("A6").HorizontalAlignment = -4108
$addedWs.Range("A6").VerticalAlignment = -4160

# --- Populate "Changed Books" sheet (original content preserved) ---
$changedWs.Range("A1").Value = "Name"
$changedWs.Range("B1").Value = "Changes"
$changedWs.Range("A1:B1").Font.Bold = $true
$changedWs.Range("A1:B1").Borders.LineStyle = 1
$changedWs.Range("A1:B1").HorizontalAlignment = -4108
$changedWs.Range("A1:B1").VerticalAlignment = -4160

$changedWs.Range("A2").Value = "Bản Đồ"
$changedWs.Range("B2").Value = "{'Price (vnd)': '224250 → 216030.76', 'Discount (%)': '35 → 23', 'Sold': '5935 → 5932'}"
$changedWs.Range("A3").Value = "Cây Cam Ngọt Của Tôi"
$changedWs.Range("B3").Value = "{'Price (vnd)': '64800 → 68046.35', 'Discount (%)': '40 → 24', 'Sold': '72191 → 72200'}"
$changedWs.Range("A4").Value = "Không Phải Sói Nhưng Cũng Đừng Là Cừu -Tặng kèm bookmark 2 mặt"
$changedWs.Range("B4").Value = "{'Price (vnd)': '85760 → 93128.02', 'Discount (%)': '33 → 14', 'Sold': '12937 → 12928'}"
$changedWs.Range("A5").Value = "Dear, Darling"
$changedWs.Range("B5").Value = "{'Price (vnd)': '58960 → 63200.76', 'Discount (%)': '33 → 2', 'Sold': '6084 → 6095'}"
$changedWs.Range("A6").Value = "Sapiens Lược Sử Loài Người (Tái Bản)"
$changedWs.Range("B6").Value = "{'Price (vnd)': '251650 → 262031.99', 'Discount (%)': '0 → 23', 'Sold': '77 → 70'}"
$changedWs.Range("A7").Value = "Một Thoáng Ta Rực Rỡ Ở Nhân Gian"
$changedWs.Range("B7").Value = "{'Price (vnd)': '81000 → 81489.13', 'Discount (%)': '40 → 2', 'Sold': '13608 → 13618'}"
$changedWs.Range("A8").Value = "Những Tù Nhân Của Địa Lý"
$changedWs.Range("B8").Value = "{'Price (vnd)': '126000 → 119782.47'}"
$changedWs.Range("A9").Value = "Đại Dương Đen - Những Câu Chuyện Từ Thế Giới Của Trầm Cảm"
$changedWs.Range("B9").Value = "{'Price (vnd)': '144000 → 149642.88', 'Discount (%)': '40 → 26', 'Sold': '12918 → 12924'}"
$changedWs.Range("A10").Value = "Rèn Luyện Tư Duy Phản Biện"
$changedWs.Range("B10").Value = "{'Price (vnd)': '59400 → 55297.33', 'Discount (%)': '40 → 29', 'Sold': '30363 → 30371'}"
$changedWs.Range("A11").Value = "Điềm Tĩnh Và Nóng Giận"
$changedWs.Range("B11").Value = "{'Price (vnd)': '64320 → 59851.86', 'Discount (%)': '33 → 19', 'Sold': '10744 → 10752'}"
$changedWs.Range("A12").Value = "Không Ai Có Thể Làm Bạn Tổn Thương Trừ Khi Bạn Cho Phép (Tặng Kèm 1 Bộ Bookmark Gồm 4 Cái)"
$changedWs.Range("B12").Value = "{'Price (vnd)': '73780 → 70043.44', 'Discount (%)': '38 → 8', 'Sold': '11904 → 11896'}"
$changedWs.Range("A13").Value = "Yêu Những Ngày Nắng Chẳng Ghét Những Ngày Mưa"
$changedWs.Range("B13").Value = "{'Price (vnd)': '80520 → 84128.91', 'Discount (%)': '39 → 28', 'Sold': '5318 → 5319'}"
$changedWs.Range("A14").Value = "Tâm Lý Học - Phác Họa Chân Dung Kẻ Phạm Tội"
$changedWs.Range("B14").Value = "{'Price (vnd)': '89900 → 97951.47', 'Discount (%)': '38 → 32', 'Sold': '22215 → 22219'}"
$changedWs.Range("A15").Value = "Xứ Cát"
$changedWs.Range("B15").Value = "{'Price (vnd)': '149400 → 158464.09', 'Discount (%)': '40 → 46', 'Sold': '7208 → 7213'}"
$changedWs.Range("A16").Value = "Đắc Nhân Tâm (Bìa Mềm)(Tái Bản)"
$changedWs.Range("B16").Value = "{'Price (vnd)': '56300 → 53496.25', 'Discount (%)': '35 → 25', 'Sold': '10621 → 10612'}"

Write-Output "done"
